# "final model with video"
# Replace the training/validation loss samples with the final (3-epoch)
# run and drop the now-unused trailing rows, updating the chart series
# to match the shrunk data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New loss values (epochs 1-3 of the final run)
$ws.Range("A2").Value = 0.0255
$ws.Range("B2").Value = 0.0215
$ws.Range("A3").Value = 0.0213
$ws.Range("B3").Value = 0.0203
$ws.Range("A4").Value = 0.0197
$ws.Range("B4").Value = 0.019

# Rows 5 and 6 are no longer part of the data set
$ws.Rows("5:6").Delete()

# Point the chart series at the new, smaller data range
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$s1 = $chart.SeriesCollection(1)
$s1.Formula = "=SERIES(Sheet1!`$A`$1,,Sheet1!`$A`$2:`$A`$4,1)"
$s2 = $chart.SeriesCollection(2)
$s2.Formula = "=SERIES(Sheet1!`$B`$1,,Sheet1!`$B`$2:`$B`$4,2)"

# Match the author's final selection
$ws.Range("D14").Select()
